$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 5.315642833709717
$ws.Range("C3").Value = 5.31564474105835
$ws.Range("C4").Value = 5.3156418800354
$ws.Range("C5").Value = 5.315646648406982
$ws.Range("C6").Value = 5.315648078918457
$ws.Range("C7").Value = 5.315647602081299
$ws.Range("C8").Value = 5.315646648406982
$ws.Range("C9").Value = 5.315647602081299
$ws.Range("C10").Value = 5.315646648406982
$ws.Range("C11").Value = 5.315646648406982
$ws.Range("C12").Value = 5.315646648406982
$ws.Range("C13").Value = 5.315649032592773
$ws.Range("C14").Value = 5.315650463104248
$ws.Range("C15").Value = 5.315649509429932
$ws.Range("C16").Value = 5.315650463104248
$ws.Range("C17").Value = 5.315648555755615
$ws.Range("C18").Value = 5.315648555755615
$ws.Range("C19").Value = 5.315650463104248
$ws.Range("C20").Value = 5.315648555755615
$ws.Range("C21").Value = 5.315647602081299
$ws.Range("C22").Value = 5.315646648406982
$ws.Range("C23").Value = 5.315648555755615
$ws.Range("C24").Value = 5.315650463104248
$ws.Range("C25").Value = 5.315648555755615
$ws.Range("C26").Value = 5.31564998626709
$ws.Range("C27").Value = 5.315650463104248
$ws.Range("C28").Value = 5.315651416778564
$ws.Range("C29").Value = 5.315650463104248
$ws.Range("C30").Value = 5.315651416778564
$ws.Range("C31").Value = 5.315651416778564
$ws.Range("C32").Value = 5.315652370452881
$ws.Range("C33").Value = 5.315651893615723
